$wb = $excel.ActiveWorkbook

# "Push choices sheet display.text into display.title.text"
# The "choices" sheet header cell C1 currently reads "display.text";
# rename it to "display.title.text" (merging with the existing shared
# string used elsewhere, e.g. on the "settings" sheet).
$choices = $wb.Worksheets.Item("choices")
$choices.Range("C1").Value = "display.title.text"

# Also push display.image into display.title.image (no such column exists
# on this particular workbook's choices sheet, so there is nothing further
# to rename here).

# The choices sheet becomes the active/selected tab in the saved workbook.
$choices.Activate()
